$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Relabel the boolean-ish "status" values with friendlier Portuguese text
# (shared strings "True"/"False" -> "ativo"/"inativo").
$ws.Cells.Replace("True", "ativo") | Out-Null
$ws.Cells.Replace("False", "inativo") | Out-Null

# Move the active selection to F4.
$ws.Range("F4").Select() | Out-Null

# Configure page setup for printing (A4 portrait), as part of enabling
# uploads of the other worksheets too.
$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1
